$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 816; this shifts the existing rows 816-860 down to 817-861
$ws.Rows(816).Insert()

# Populate the newly inserted row 816 with the new weekly price record
$ws.Cells.Item(816, 1).Value2 = 6
$ws.Cells.Item(816, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(816, 3).Value2 = "Metropolitana"
$ws.Cells.Item(816, 4).Value2 = 45267
$ws.Cells.Item(816, 5).Value2 = 13
$ws.Cells.Item(816, 6).Value2 = 100112012
$ws.Cells.Item(816, 7).Value2 = "Espinaca"
$ws.Cells.Item(816, 8).Value2 = "Sin especificar"
$ws.Cells.Item(816, 9).Value2 = "Primera"
$ws.Cells.Item(816, 10).Value2 = 450
$ws.Cells.Item(816, 11).Value2 = 6500
$ws.Cells.Item(816, 12).Value2 = 7000
$ws.Cells.Item(816, 13).Value2 = 6778
$ws.Cells.Item(816, 14).Value2 = "$/cuna 10 kilos"
$ws.Cells.Item(816, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(816, 16).Value2 = 678
$ws.Cells.Item(816, 17).Value2 = 10
$ws.Cells.Item(816, 18).Value2 = "Hortaliza"
